$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in column H1, matching the style used by the other
# header cells (copy format from G1, the "sum" header).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill H2:H8 with 0 (new "Save" column values)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
